# Session2.pptx edit: rename presenter "Bea" -> "Leonardo" and merge the
# "March" / " 2024" date runs into a single "March 2024" run, on both the
# Slide 1 (welcome) and Slide 2 (session) title-card shapes ("CustomShape 2").

$p = $ppt.ActivePresentation

function Update-PresenterCard($slideIndex) {
    $slide = $p.Slides.Item($slideIndex)
    $shape = $slide.Shapes.Item(2)
    $tf = $shape.TextFrame
    $tr = $tf.TextRange

    # Paragraph 1: "Irina & Bea" -> "Irina & Leonardo"
    $para1 = $tr.Paragraphs(1)
    $para1.Runs(2).Text = "Leonardo"

    # Paragraph 2: "March" + " 2024" (two runs) -> single run "March 2024"
    # carrying the second run's character formatting (b="0" strike="noStrike").
    $para2 = $tr.Paragraphs(2)
    # Force a genuine text replacement (no shared prefix/suffix with the old
    # two-run text) so the engine collapses the paragraph down to one run.
    $para2.Text = "##TMP##"
    $para2.Text = "March 2024"
    $para2.Runs(1).Font.Bold = 0
    $para2.Runs(1).Font.Strikethrough = 0
}

Update-PresenterCard 1
Update-PresenterCard 2
